$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-7: sending-cluster column (A) and the numeric NATMI metrics
# that were recomputed against the new TPM values.
$ws.Range("A2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.9795656666666667
$ws.Range("H2").Value = 2.938697
$ws.Range("I2").Value = 0.2969748092101394
$ws.Range("J2").Value = 0.2969748092101394
$ws.Range("M2").Value = 91.51130433333333
$ws.Range("N2").Value = 274.533913
$ws.Range("O2").Value = 0.9685519820468944
$ws.Range("P2").Value = 0.9685519820468945
$ws.Range("Q2").Value = 89.64133183681788
$ws.Range("R2").Value = 806.7719865313609
$ws.Range("S2").Value = 0.2876355400784788
$ws.Range("T2").Value = 0.2876355400784789
$ws.Range("A3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.9795656666666667
$ws.Range("H3").Value = 2.938697
$ws.Range("I3").Value = 0.2969748092101394
$ws.Range("J3").Value = 0.2969748092101394
$ws.Range("O3").Value = 0.001425786415744213
$ws.Range("P3").Value = 0.001425786415744214
$ws.Range("Q3").Value = 0.131959250088
$ws.Range("R3").Value = 1.187633250792
$ws.Range("S3").Value = 0.0004234226487900462
$ws.Range("T3").Value = 0.0004234226487900464
$ws.Range("A4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.9795656666666667
$ws.Range("H4").Value = 2.938697
$ws.Range("I4").Value = 0.2969748092101394
$ws.Range("J4").Value = 0.2969748092101394
$ws.Range("M4").Value = 2.836578333333333
$ws.Range("N4").Value = 8.509734999999999
$ws.Range("O4").Value = 0.03002223153736139
$ws.Range("P4").Value = 0.03002223153736139
$ws.Range("Q4").Value = 2.778614746143889
$ws.Range("R4").Value = 25.007532715295
$ws.Range("S4").Value = 0.00891584648287053
$ws.Range("T4").Value = 0.00891584648287053
$ws.Range("A5").Value = "MuSCs"
$ws.Range("G5").Value = 2.318915
$ws.Range("H5").Value = 6.956745
$ws.Range("I5").Value = 0.7030251907898606
$ws.Range("J5").Value = 0.7030251907898607
$ws.Range("M5").Value = 91.51130433333333
$ws.Range("N5").Value = 274.533913
$ws.Range("O5").Value = 0.9685519820468944
$ws.Range("P5").Value = 0.9685519820468945
$ws.Range("Q5").Value = 212.2069362881317
$ws.Range("R5").Value = 1909.862426593185
$ws.Range("S5").Value = 0.6809164419684156
$ws.Range("T5").Value = 0.6809164419684157
$ws.Range("A6").Value = "MuSCs"
$ws.Range("G6").Value = 2.318915
$ws.Range("H6").Value = 6.956745
$ws.Range("I6").Value = 0.7030251907898606
$ws.Range("J6").Value = 0.7030251907898607
$ws.Range("O6").Value = 0.001425786415744213
$ws.Range("P6").Value = 0.001425786415744214
$ws.Range("Q6").Value = 0.31238567748
$ws.Range("R6").Value = 2.81147109732
$ws.Range("S6").Value = 0.001002363766954167
$ws.Range("T6").Value = 0.001002363766954168
$ws.Range("A7").Value = "MuSCs"
$ws.Range("G7").Value = 2.318915
$ws.Range("H7").Value = 6.956745
$ws.Range("I7").Value = 0.7030251907898606
$ws.Range("J7").Value = 0.7030251907898607
$ws.Range("M7").Value = 2.836578333333333
$ws.Range("N7").Value = 8.509734999999999
$ws.Range("O7").Value = 0.03002223153736139
$ws.Range("P7").Value = 0.03002223153736139
$ws.Range("Q7").Value = 6.577784045841667
$ws.Range("R7").Value = 59.20005641257499
$ws.Range("S7").Value = 0.02110638505449087
$ws.Range("T7").Value = 0.02110638505449087
# Rows 8-10 (the old "MuSCs" sending-cluster block) are no longer present
# in the updated output -- remove them so the sheet ends at row 7.
$ws.Range("A8:T10").Delete()
